$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.423.94'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.839.31'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.65%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +2.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.45'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.45%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4360'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.79%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3721'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.51%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8722'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.42'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.969.76'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +10.90%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.682'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07155'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.59%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.030'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008998'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.024'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.36%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.451.93'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.253'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.16'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.162.71'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +8.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.14'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.899'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.258'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.921'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.34'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09038'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.200'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +7.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7588'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.469'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.12%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.09%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.150'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01957'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05253'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5165'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.797'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.50%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.528'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.465'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +5.77%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'PaxosStandard'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.027'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.22%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '108.73'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.65%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.48'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.16%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.027'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.43%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4632'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.34%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.670'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06292'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.49%  '
